# edit.ps1 -- apply the "Device login with exceptios and logs" change set
#
# Summary of edits (see accompanying diff):
#   1. Wrap the misspelled "Setps" heading run in spell-check proofErr markers.
#   2. Append " ()" after "Have SQL error E104".
#   3. Append " (Not have Authentication)" after "...authorized E103".
#   4. Insert eight new list-item paragraphs (same list level as the SQL-error
#      bullets) describing further error/exception codes, the second of which
#      wraps the misspelled "priviledge" in spell-check proofErr markers and
#      the fourth of which carries a <w:lastRenderedPageBreak/>. The very last
#      new bullet ("telephone ... E108") lands in the paragraph that already
#      holds the _GoBack bookmark, exactly as in the target document.
#
# Because the Word object model has no direct "add a proofErr marker" verb,
# we drive every paragraph body through Range.InsertXML using the
# WordprocessingML "Word XML Package" envelope -- this replaces the full
# contents of the paragraph (whatever pPr/runs it starts with) with exactly
# the OOXML fragment we specify, so we always pass the paragraph's own pPr
# back through unchanged alongside whatever new runs/proofErr markers we need.
#
# Note: Paragraph.Range.Text includes the trailing paragraph-mark character,
# so exact "-eq" string comparisons never match -- "-like "*...*"" is used
# throughout when locating paragraphs by their visible text.

$d = $word.ActiveDocument

function New-ParaXml([string]$innerBodyXml) {
    return "<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>$innerBodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

$lvl1PPr = "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr>"

# ---------------------------------------------------------------------
# 1) "Setps" -> wrap in proofErr spellStart/spellEnd (ilvl 0 heading)
# ---------------------------------------------------------------------
$setpsPara = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Setps*") {
        $setpsPara = $i
        break
    }
}
$r = $d.Paragraphs($setpsPara).Range
$xml = New-ParaXml("<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:proofErr w:type='spellStart'/><w:r><w:t>Setps</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>")
$r.InsertXML($xml)

# ---------------------------------------------------------------------
# 2) "Have SQL error E104" -> append " ()"
# ---------------------------------------------------------------------
$sqlPara = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Have SQL error E104*") {
        $sqlPara = $i
        break
    }
}
$r = $d.Paragraphs($sqlPara).Range
$xml = New-ParaXml("<w:p>$lvl1PPr<w:r><w:t>Have SQL error E104</w:t></w:r><w:r><w:t xml:space='preserve'> ()</w:t></w:r></w:p>")
$r.InsertXML($xml)

# ---------------------------------------------------------------------
# 3) "... authorized E103" -> append " (Not have Authentication)"
#    (this is the paragraph right after the SQL-error one)
# ---------------------------------------------------------------------
$authPara = $sqlPara + 1
$r = $d.Paragraphs($authPara).Range
$xml = New-ParaXml("<w:p>$lvl1PPr<w:r><w:t xml:space='preserve'> </w:t></w:r><w:r><w:t xml:space='preserve'>device is not </w:t></w:r><w:r><w:t>authorized E103</w:t></w:r><w:r><w:t xml:space='preserve'> (</w:t></w:r><w:r><w:t>Not have Authentication</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>")
$r.InsertXML($xml)

# ---------------------------------------------------------------------
# 4) Insert the first seven new bullet paragraphs right before the
#    paragraph that holds the _GoBack bookmark (the paragraph right
#    after "... authorized E103", which is currently empty).
# ---------------------------------------------------------------------
$goBackPara = $authPara + 1

$newBulletBodies = @(
    "<w:r><w:t xml:space='preserve'>No </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>priviledge</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> </w:t></w:r><w:r><w:t xml:space='preserve'> </w:t></w:r><w:r><w:t>E101 (</w:t></w:r><w:r><w:t>Not have Authentication</w:t></w:r><w:r><w:t>)</w:t></w:r>",
    "<w:r><w:lastRenderedPageBreak/><w:t>Error in background image E102</w:t></w:r>",
    "<w:r><w:t xml:space='preserve'>Data you send may have error </w:t></w:r><w:r><w:t>E100</w:t></w:r>",
    "<w:r><w:t>Login failed E105</w:t></w:r>",
    "<w:r><w:t>User is blocked E106</w:t></w:r>",
    "<w:r><w:t>User want to be Active E107</w:t></w:r>",
    "<w:r><w:t>Email</w:t></w:r><w:r><w:t xml:space='preserve'> want to be Active E108</w:t></w:r>"
)

foreach ($body in $newBulletBodies) {
    $insertionPoint = $d.Paragraphs($goBackPara).Range
    $insertionPoint.Collapse(1)
    $insertionPoint.InsertParagraphBefore()
    $newParaIndex = $goBackPara
    $goBackPara = $goBackPara + 1
    $newRange = $d.Paragraphs($newParaIndex).Range
    $xml = New-ParaXml("<w:p>$lvl1PPr$body</w:p>")
    $newRange.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 5) Final bullet ("telephone ... E108") goes into the paragraph that
#    already holds the _GoBack bookmark -- add the two runs, keep the
#    bookmark markers exactly where they were (after the new text).
# ---------------------------------------------------------------------
$finalPara = $d.Paragraphs($goBackPara).Range
$xml = New-ParaXml("<w:p>$lvl1PPr<w:r><w:t>telephone</w:t></w:r><w:r><w:t xml:space='preserve'> want to be Active E108</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>")
$finalPara.InsertXML($xml)
